$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original column A (TAXON count values styled as headers) is removed;
# every remaining column shifts one place to the left (B->A, C->B, D->C, E->D, F->E).
$ws.Columns("A").Delete()
